$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to match the workbook's title.
$ws.Name = "EL at 50Apercm2"

# The wavelength column (A) and the two "Osram" intensity columns (C, D)
# were widened to fit their numeric contents.
$ws.Columns("A:A").AutoFit()
$ws.Columns("C:C").AutoFit()
$ws.Columns("D:D").AutoFit()

# Leave the active selection where the author last clicked.
$ws.Range("I13").Select()
